$d = $word.ActiveDocument

# Locate the (unique) run that currently reads:
#   "co 3 spotkanie niebezpieczeństwa tracą 2 członków populacji"
$rng = $d.Content
$found = $rng.Find.Execute("tracą 2")
if (-not $found) {
    throw "Could not find anchor text 'tracą 2'"
}

# Grab the formatting (blue text, same size) already used by that run so the
# inserted words keep matching look-and-feel.
$clr  = $rng.Font.Color
$size = $rng.Font.Size

# Collapse to the point right after "...tracą 2" (before " członków populacji")
$rng.Collapse(0)
$insStart = $rng.Start

$insertion = " razy wiecej"
$razyText   = " razy "
$wiecejText = "wiecej"

# Type the new words in, same as a user would.
$rng.InsertAfter($insertion)

$razyStart   = $insStart
$wiecejStart = $insStart + $razyText.Length
$wiecejEnd   = $wiecejStart + $wiecejText.Length

# Find the tail run (" członków populacji") that follows what we just typed,
# so we can re-stamp its formatting too and keep it as a distinct run.
$tailRng = $d.Range($wiecejEnd, $d.Content.End)
$tailFound = $tailRng.Find.Execute(" członków populacji")
if (-not $tailFound) {
    throw "Could not find tail text ' członków populacji'"
}
$tailStart = $tailRng.Start
$tailEnd = $tailRng.End

# Re-apply the existing formatting explicitly on each of the three newly
# delimited spans (the text typed plus the pre-existing tail), nudging each
# one into its own run - mirroring how Word keeps a freshly
# typed/autocorrect-flagged word such as "wiecej" in its own run, separate
# from its neighbours - while leaving the visible formatting unchanged.
$razySeg = $d.Range($razyStart, $wiecejStart)
$razySeg.Font.Bold = $true
$razySeg.Font.Bold = $false
$razySeg.Font.Color = $clr
$razySeg.Font.Size = $size

$wiecejSeg = $d.Range($wiecejStart, $wiecejEnd)
$wiecejSeg.Font.Bold = $true
$wiecejSeg.Font.Bold = $false
$wiecejSeg.Font.Color = $clr
$wiecejSeg.Font.Size = $size

$tailSeg = $d.Range($tailStart, $tailEnd)
$tailSeg.Font.Bold = $true
$tailSeg.Font.Bold = $false
$tailSeg.Font.Color = $clr
$tailSeg.Font.Size = $size

Write-Output "Inserted '$insertion' after 'tracą 2'"
